$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the (soon to be
#    deleted) "git push --tags" paragraph to the end of the
#    "...: define STAGING_BACKEND." paragraph.
#
#    The bookmark must end up *after* the existing last run of that
#    paragraph (not splitting it), i.e. immediately before the
#    paragraph's closing tag.  Adding a bookmark with a Range whose
#    Start sits exactly on the paragraph-mark character is mishandled
#    by this runtime, so we work around it by briefly inserting a
#    placeholder character after the paragraph's text, wrapping that
#    placeholder with the bookmark, and then deleting the placeholder
#    again - this leaves the now-collapsed bookmark in exactly the
#    right spot without disturbing the existing run.
# ------------------------------------------------------------------
$backendParagraph = $d.Paragraphs.Item(3)
$endPos = $backendParagraph.Range.End
$insertionPoint = $d.Range($endPos - 1, $endPos - 1)
$insertionPoint.InsertAfter("X")
$placeholder = $d.Range($endPos - 1, $endPos)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder = $d.Range($endPos - 1, $endPos)
$placeholder.Delete()

# ------------------------------------------------------------------
# 2) Remove the "and CloudSetupApiSyncSampleSupport" reference from
#    the "Change the CloudApiPublic ... Assembly and File Versions"
#    paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute(" and CloudSetupApiSyncSampleSupport ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Delete whole paragraphs that are no longer part of the
#    procedure:
#      - "Change the CloudSetupSdkSyncSample version. ..."
#      - "Exit Visual Studio."
#      - "To delete a tag locally and remotely:"
#      - "git tag -d 12345"
#      - "git push origin :refs/tags/12345"
#      - "To create a local tag named 12345 and push it to remote:"
#      - "git tag 12345"
#      - "git push --tags"
#    (deleted from the bottom up so earlier paragraph indices stay
#    valid while iterating)
# ------------------------------------------------------------------
$paragraphsToDelete = @(15, 14, 13, 12, 11, 10, 7, 5)
foreach ($idx in $paragraphsToDelete) {
    $p = $d.Paragraphs.Item($idx)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.Delete()
}

# ------------------------------------------------------------------
# 4) Reword "Commit and push to master again." to
#    "If anything needed to be changed, commit and push to master
#    again."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Commit and push to master again.", $true, $false, $false, $false, $false, $true, 1, $false, "If anything needed to be changed, commit and push to master again.", 2) | Out-Null
